# Apply updated currentAveragePrice / Leve profit figures from the latest
# scheduled market-data refresh (Sheets via scheduled runner).
$wb = $excel.ActiveWorkbook

# ================= Sheet ALC =================
$ws = $wb.Sheets.Item("ALC")
# Row 113
$ws.Range("H113").Value = 3262.842
$ws.Range("I113").Value = 3294.4443
$ws.Range("J113").Value = 3234.4
$ws.Range("K113").Value = 3294.4443
$ws.Range("L113").Value = 3234.4
$ws.Range("M113").Value = -40.44430000000011
$ws.Range("N113").Value = -9742.4
# Row 116
$ws.Range("H116").Value = 4085.625
$ws.Range("I116").Value = 3877.2222
$ws.Range("J116").Value = 4210.6665
$ws.Range("K116").Value = 3877.2222
$ws.Range("L116").Value = 4210.6665
$ws.Range("M116").Value = -435.2222000000002
$ws.Range("N116").Value = -11094.6665
# Row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
# Row 132
$ws.Range("H132").Value = 3848463.8
$ws.Range("I132").Value = 4083841.2
$ws.Range("J132").Value = 3966.6667
$ws.Range("K132").Value = 12251523.6
$ws.Range("L132").Value = 11900.0001
$ws.Range("M132").Value = -12248993.6
$ws.Range("N132").Value = -16960.0001
# Row 138
$ws.Range("H138").Value = 4315.0435
$ws.Range("J138").Value = 4870.3945
$ws.Range("L138").Value = 14611.1835
$ws.Range("N138").Value = -24891.1835

# ================= Sheet ARM =================
$ws = $wb.Sheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 4123.3335
$ws.Range("J3").Value = 1500
$ws.Range("L3").Value = 1500
$ws.Range("N3").Value = -1730
# Row 74
$ws.Range("H74").Value = 1656.091
$ws.Range("I74").Value = 1332.7778
$ws.Range("J74").Value = 3111
$ws.Range("K74").Value = 1332.7778
$ws.Range("L74").Value = 3111
$ws.Range("M74").Value = -458.7778000000001
$ws.Range("N74").Value = -4859
# Row 77
$ws.Range("H77").Value = 1656.091
$ws.Range("I77").Value = 1332.7778
$ws.Range("J77").Value = 3111
$ws.Range("K77").Value = 6663.889
$ws.Range("L77").Value = 15555
$ws.Range("M77").Value = -2295.889
$ws.Range("N77").Value = -24291
# Row 110
$ws.Range("H110").Value = 1543.1333
$ws.Range("I110").Value = 603.4583
$ws.Range("J110").Value = 5301.8335
$ws.Range("K110").Value = 603.4583
$ws.Range("L110").Value = 5301.8335
$ws.Range("M110").Value = 1441.5417
$ws.Range("N110").Value = -9391.833500000001

# ================= Sheet BSM =================
$ws = $wb.Sheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1740.1052
$ws.Range("I105").Value = 1406.1538
$ws.Range("J105").Value = 2463.6667
$ws.Range("K105").Value = 1406.1538
$ws.Range("L105").Value = 2463.6667
$ws.Range("M105").Value = 340.8462
$ws.Range("N105").Value = -5957.6667

# ================= Sheet CRP =================
$ws = $wb.Sheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 64002
$ws.Range("J4").Value = 64002
$ws.Range("L4").Value = 64002
$ws.Range("N4").Value = -64226
# Row 31
$ws.Range("H31").Value = 2661.9773
$ws.Range("I31").Value = 1794.1455
$ws.Range("J31").Value = 4108.364
$ws.Range("K31").Value = 1794.1455
$ws.Range("L31").Value = 4108.364
$ws.Range("M31").Value = -1499.1455
$ws.Range("N31").Value = -4698.364
# Row 34
$ws.Range("H34").Value = 2661.9773
$ws.Range("I34").Value = 1794.1455
$ws.Range("J34").Value = 4108.364
$ws.Range("K34").Value = 1794.1455
$ws.Range("L34").Value = 4108.364
$ws.Range("M34").Value = -1592.1455
$ws.Range("N34").Value = -4512.364
# Row 107
$ws.Range("H107").Value = 1305.04
$ws.Range("I107").Value = 1098.3125
$ws.Range("J107").Value = 1672.5555
$ws.Range("K107").Value = 1098.3125
$ws.Range("L107").Value = 1672.5555
$ws.Range("M107").Value = 821.6875
$ws.Range("N107").Value = -5512.5555
# Row 132
$ws.Range("H132").Value = 2567.2727
$ws.Range("I132").Value = 2269.68
$ws.Range("J132").Value = 3497.25
$ws.Range("K132").Value = 6809.039999999999
$ws.Range("L132").Value = 10491.75
$ws.Range("M132").Value = -4279.039999999999
$ws.Range("N132").Value = -15551.75

# ================= Sheet CUL =================
$ws = $wb.Sheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 7199.0835
$ws.Range("J3").Value = 19000
$ws.Range("L3").Value = 57000
$ws.Range("N3").Value = -57224
# Row 87
$ws.Range("H87").Value = 9841.6
$ws.Range("I87").Value = 4100
$ws.Range("J87").Value = 15583.2
$ws.Range("K87").Value = 12300
$ws.Range("L87").Value = 46749.60000000001
$ws.Range("M87").Value = -11052
$ws.Range("N87").Value = -49245.60000000001
# Row 90
$ws.Range("H90").Value = 9841.6
$ws.Range("I90").Value = 4100
$ws.Range("J90").Value = 15583.2
$ws.Range("K90").Value = 36900
$ws.Range("L90").Value = 140248.8
$ws.Range("M90").Value = -30660
$ws.Range("N90").Value = -152728.8
# Row 105
$ws.Range("H105").Value = 4000
$ws.Range("J105").Value = 4000
$ws.Range("L105").Value = 12000
$ws.Range("N105").Value = -17242
# Row 125
$ws.Range("H125").Value = 2321.0667
$ws.Range("I125").Value = 1538.3334
$ws.Range("J125").Value = 2842.889
$ws.Range("K125").Value = 4615.0002
$ws.Range("L125").Value = 8528.667000000001
$ws.Range("M125").Value = 304.9997999999996
$ws.Range("N125").Value = -18368.667
# Row 131
$ws.Range("H131").Value = 1638.8334
$ws.Range("I131").Value = 3625.7144
$ws.Range("J131").Value = 1241.4572
$ws.Range("K131").Value = 10877.1432
$ws.Range("L131").Value = 3724.3716
$ws.Range("M131").Value = -5837.143199999999
$ws.Range("N131").Value = -13804.3716
# Row 133
$ws.Range("H133").Value = 6405.4546

# ================= Sheet GSM =================
$ws = $wb.Sheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 83904
$ws.Range("I4").Value = 3404
$ws.Range("J4").Value = 100004
$ws.Range("K4").Value = 3404
$ws.Range("L4").Value = 100004
$ws.Range("M4").Value = -3292
$ws.Range("N4").Value = -100228
# Row 70
$ws.Range("H70").Value = 4560.3706
$ws.Range("I70").Value = 4500
$ws.Range("J70").Value = 4681.1113
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 4681.1113
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -5221.1113
# Row 73
$ws.Range("H73").Value = 4560.3706
$ws.Range("I73").Value = 4500
$ws.Range("J73").Value = 4681.1113
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 4681.1113
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -6553.1113

# ================= Sheet LTW =================
$ws = $wb.Sheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2670.244
$ws.Range("I132").Value = 1936.7084
$ws.Range("J132").Value = 3705.8235
$ws.Range("K132").Value = 5810.1252
$ws.Range("L132").Value = 11117.4705
$ws.Range("M132").Value = -3280.1252
$ws.Range("N132").Value = -16177.4705

# ================= Sheet WVR =================
$ws = $wb.Sheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 848.6
$ws.Range("I113").Value = 176.57143
$ws.Range("J113").Value = 1436.625
$ws.Range("K113").Value = 529.71429
$ws.Range("L113").Value = 4309.875
$ws.Range("M113").Value = 1640.28571
$ws.Range("N113").Value = -8649.875
# Row 126
$ws.Range("H126").Value = 27734.9
$ws.Range("I126").Value = 49881.57
$ws.Range("J126").Value = 3257
$ws.Range("K126").Value = 149644.71
$ws.Range("L126").Value = 9771
$ws.Range("M126").Value = -147174.71
$ws.Range("N126").Value = -14711
